$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark newly completed backlog items as "Done" (column I) and record their Sprint (column J)
$ws.Range("I24").Value = "Done"
$ws.Range("I24").Interior.Color = 5287936
$ws.Range("J24").Value = "Sprint 3"

$ws.Range("I25").Value = "Done"
$ws.Range("I25").Interior.Color = 5287936
$ws.Range("J25").Value = "Sprint 2"

$ws.Range("I26").Value = "Done"
$ws.Range("I26").Interior.Color = 5287936

$ws.Range("I27").Value = "Done"
$ws.Range("I27").Interior.Color = 5287936
$ws.Range("J27").Value = "Sprint 2"

$ws.Range("I28").Value = "Done"
$ws.Range("I28").Interior.Color = 5287936
$ws.Range("J28").Value = "Sprint 2"

$ws.Range("I29").Value = "Done"
$ws.Range("I29").Interior.Color = 5287936
$ws.Range("J29").Value = "Sprint 2"

$ws.Range("I30").Value = "Done"
$ws.Range("I30").Interior.Color = 5287936
$ws.Range("J30").Value = "Sprint 2"

$ws.Range("I31").Value = "Done"
$ws.Range("I31").Interior.Color = 5287936
$ws.Range("J31").Value = "Sprint 2"

$ws.Range("I32").Value = "Done"
$ws.Range("I32").Interior.Color = 5287936
$ws.Range("J32").Value = "Sprint 2"

$ws.Range("I33").Value = "Done"
$ws.Range("I33").Interior.Color = 5287936
$ws.Range("J33").Value = "Sprint 3"

$ws.Range("I34").Value = "Done"
$ws.Range("I34").Interior.Color = 5287936
$ws.Range("J34").Value = "Sprint 3"

$ws.Range("I35").Value = "Done"
$ws.Range("I35").Interior.Color = 5287936
$ws.Range("J35").Value = "Sprint 3"

$ws.Range("I36").Value = "Done"
$ws.Range("I36").Interior.Color = 5287936
$ws.Range("J36").Value = "Sprint 2"

$ws.Range("I37").Value = "Done"
$ws.Range("I37").Interior.Color = 5287936
$ws.Range("J37").Value = "Sprint 3"

$ws.Range("I38").Value = "Done"
$ws.Range("I38").Interior.Color = 5287936
$ws.Range("J38").Value = "Sprint 2"

$ws.Range("I39").Value = "Done"
$ws.Range("I39").Interior.Color = 5287936
$ws.Range("J39").Value = "Sprint 2"

$ws.Range("I42").Value = "Done"
$ws.Range("I42").Interior.Color = 5287936
$ws.Range("J42").Value = "Sprint 2"

$ws.Range("I43").Value = "Done"
$ws.Range("I43").Interior.Color = 5287936
$ws.Range("J43").Value = "Sprint 2"

# Adjust row heights to match re-flowed content after the edits above
$ws.Rows.Item(21).RowHeight = 114
$ws.Rows.Item(24).RowHeight = 71.25
$ws.Rows.Item(25).RowHeight = 85.5
$ws.Rows.Item(43).RowHeight = 71.25
$ws.Rows.Item(45).RowHeight = 71.25
$ws.Rows.Item(46).RowHeight = 85.5
$ws.Rows.Item(47).RowHeight = 71.25
$ws.Rows.Item(51).RowHeight = 114
$ws.Rows.Item(54).RowHeight = 71.25
$ws.Rows.Item(57).RowHeight = 71.25

